$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 67 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2451.2273
$ws.Range("I15").Value = 2451.2273
$ws.Range("K15").Value = 7353.6819
$ws.Range("M15").Value = -7184.6819
$ws.Range("H49").Value = 399.33334
$ws.Range("I49").Value = 399.33334
$ws.Range("K49").Value = 1198.00002
$ws.Range("M49").Value = -1062.00002
$ws.Range("H62").Value = 8337836.5
$ws.Range("I62").Value = 13891074
$ws.Range("K62").Value = 13891074
$ws.Range("M62").Value = -13890450
$ws.Range("H65").Value = 8337836.5
$ws.Range("I65").Value = 13891074
$ws.Range("K65").Value = 69455370
$ws.Range("M65").Value = -69452250
$ws.Range("H86").Value = 3100581
$ws.Range("I86").Value = 5399
$ws.Range("J86").Value = 4390240.5
$ws.Range("K86").Value = 5399
$ws.Range("L86").Value = 4390240.5
$ws.Range("M86").Value = -4276
$ws.Range("N86").Value = -4392486.5
$ws.Range("H89").Value = 3100581
$ws.Range("I89").Value = 5399
$ws.Range("J89").Value = 4390240.5
$ws.Range("K89").Value = 26995
$ws.Range("L89").Value = 21951202.5
$ws.Range("M89").Value = -21379
$ws.Range("N89").Value = -21962434.5
$ws.Range("H100").Value = 2950
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -359
$ws.Range("N100").Value = -6082
$ws.Range("H106").Value = 400
$ws.Range("I106").Value = 400
$ws.Range("K106").Value = 400
$ws.Range("M106").Value = 231
$ws.Range("H112").Value = 3237.1667
$ws.Range("J112").Value = 3237.1667
$ws.Range("L112").Value = 9711.500100000001
$ws.Range("N112").Value = -11927.5001
$ws.Range("H113").Value = 5999.5
$ws.Range("I113").Value = 5999
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 5999
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -2745
$ws.Range("N113").Value = -12508
$ws.Range("H125").Value = 7410297.5
$ws.Range("I125").Value = 2154.5
$ws.Range("J125").Value = 12349060
$ws.Range("K125").Value = 19390.5
$ws.Range("L125").Value = 111141540
$ws.Range("M125").Value = -16930.5
$ws.Range("N125").Value = -111146460
$ws.Range("H132").Value = 7756.4443
$ws.Range("I132").Value = 10762
$ws.Range("K132").Value = 32286
$ws.Range("M132").Value = -29756
$ws.Range("H135").Value = 557101.25
$ws.Range("I135").Value = 1001419.4
$ws.Range("K135").Value = 9012774.6
$ws.Range("M135").Value = -9010239.6

# --- Sheet ARM: 45 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 390
$ws.Range("I17").Value = 390
$ws.Range("K17").Value = 390
$ws.Range("M17").Value = -217
$ws.Range("H32").Value = 4655.375
$ws.Range("I32").Value = 4103.655
$ws.Range("K32").Value = 4103.655
$ws.Range("M32").Value = -3816.655
$ws.Range("H45").Value = 2937.3845
$ws.Range("I45").Value = 1841.4286
$ws.Range("K45").Value = 1841.4286
$ws.Range("M45").Value = -1464.4286
$ws.Range("H61").Value = 4379
$ws.Range("I61").Value = 2369.8572
$ws.Range("J61").Value = 6723
$ws.Range("K61").Value = 2369.8572
$ws.Range("L61").Value = 6723
$ws.Range("M61").Value = -2157.8572
$ws.Range("N61").Value = -7147
$ws.Range("H74").Value = 1641.25
$ws.Range("I74").Value = 1457.8572
$ws.Range("K74").Value = 1457.8572
$ws.Range("M74").Value = -583.8571999999999
$ws.Range("H77").Value = 1641.25
$ws.Range("I77").Value = 1457.8572
$ws.Range("K77").Value = 7289.286
$ws.Range("M77").Value = -2921.286
$ws.Range("H102").Value = 1133
$ws.Range("I102").Value = 1099.5
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1099.5
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 522.5
$ws.Range("N102").Value = -4444
$ws.Range("H132").Value = 4637.9067
$ws.Range("I132").Value = 3656.1
$ws.Range("K132").Value = 10968.3
$ws.Range("M132").Value = -8438.299999999999
$ws.Range("H136").Value = 4379
$ws.Range("I136").Value = 2369.8572
$ws.Range("J136").Value = 6723
$ws.Range("K136").Value = 7109.571599999999
$ws.Range("L136").Value = 20169
$ws.Range("M136").Value = -4559.571599999999
$ws.Range("N136").Value = -25269

# --- Sheet BSM: 30 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2706.647
$ws.Range("I20").Value = 1955.3
$ws.Range("J20").Value = 3780
$ws.Range("K20").Value = 1955.3
$ws.Range("L20").Value = 3780
$ws.Range("M20").Value = -1708.3
$ws.Range("N20").Value = -4274
$ws.Range("H30").Value = 640
$ws.Range("I30").Value = 640
$ws.Range("K30").Value = 640
$ws.Range("M30").Value = -515
$ws.Range("H50").Value = 80777.336
$ws.Range("J50").Value = 80777.336
$ws.Range("L50").Value = 80777.336
$ws.Range("N50").Value = -81925.336
$ws.Range("H99").Value = 2339.6667
$ws.Range("I99").Value = 1457.2858
$ws.Range("J99").Value = 5428
$ws.Range("K99").Value = 1457.2858
$ws.Range("L99").Value = 5428
$ws.Range("M99").Value = 40.71419999999989
$ws.Range("N99").Value = -8424
$ws.Range("H107").Value = 371221.28
$ws.Range("I107").Value = 968.3333
$ws.Range("K107").Value = 968.3333
$ws.Range("M107").Value = 951.6667
$ws.Range("H134").Value = 39410.965
$ws.Range("I134").Value = 2585.111
$ws.Range("K134").Value = 7755.333
$ws.Range("M134").Value = -5220.333

# --- Sheet CRP: 25 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23835.447
$ws.Range("I31").Value = 1489.125
$ws.Range("J31").Value = 71507.60000000001
$ws.Range("K31").Value = 1489.125
$ws.Range("L31").Value = 71507.60000000001
$ws.Range("M31").Value = -1194.125
$ws.Range("N31").Value = -72097.60000000001
$ws.Range("H34").Value = 23835.447
$ws.Range("I34").Value = 1489.125
$ws.Range("J34").Value = 71507.60000000001
$ws.Range("K34").Value = 1489.125
$ws.Range("L34").Value = 71507.60000000001
$ws.Range("M34").Value = -1287.125
$ws.Range("N34").Value = -71911.60000000001
$ws.Range("H132").Value = 4989.154
$ws.Range("I132").Value = 4830.0586
$ws.Range("J132").Value = 5289.6665
$ws.Range("K132").Value = 14490.1758
$ws.Range("L132").Value = 15868.9995
$ws.Range("M132").Value = -11960.1758
$ws.Range("N132").Value = -20928.9995
$ws.Range("H134").Value = 913833.5600000001
$ws.Range("I134").Value = 560407
$ws.Range("K134").Value = 1681221
$ws.Range("M134").Value = -1678686

# --- Sheet CUL: 19 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 140612
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 140612
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 421836
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -422004
$ws.Range("H82").Value = 7395
$ws.Range("I82").Value = 10000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29594
$ws.Range("H85").Value = 7395
$ws.Range("I85").Value = 10000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28596
$ws.Range("H117").Value = 946.2
$ws.Range("J117").Value = 1315.5
$ws.Range("L117").Value = 3946.5
$ws.Range("N117").Value = -10830.5

# --- Sheet GSM: 25 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5987.273
$ws.Range("I80").Value = 4928.8887
$ws.Range("J80").Value = 10750
$ws.Range("K80").Value = 4928.8887
$ws.Range("L80").Value = 10750
$ws.Range("M80").Value = -3930.8887
$ws.Range("N80").Value = -12746
$ws.Range("H83").Value = 5987.273
$ws.Range("I83").Value = 4928.8887
$ws.Range("J83").Value = 10750
$ws.Range("K83").Value = 24644.4435
$ws.Range("L83").Value = 53750
$ws.Range("M83").Value = -19652.4435
$ws.Range("N83").Value = -63734
$ws.Range("H122").Value = 3529.9443
$ws.Range("I122").Value = 2141.6667
$ws.Range("J122").Value = 4918.222
$ws.Range("K122").Value = 6425.000100000001
$ws.Range("L122").Value = 14754.666
$ws.Range("M122").Value = -3975.000100000001
$ws.Range("N122").Value = -19654.666
$ws.Range("H132").Value = 369902.44
$ws.Range("I132").Value = 503653.66
$ws.Range("K132").Value = 1510960.98
$ws.Range("M132").Value = -1508430.98

# --- Sheet LTW: 25 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 867.7083
$ws.Range("I55").Value = 213.47368
$ws.Range("J55").Value = 3353.8
$ws.Range("K55").Value = 213.47368
$ws.Range("L55").Value = 3353.8
$ws.Range("M55").Value = -40.47368
$ws.Range("N55").Value = -3699.8
$ws.Range("H68").Value = 5937.25
$ws.Range("I68").Value = 4499.5
$ws.Range("J68").Value = 7375
$ws.Range("K68").Value = 4499.5
$ws.Range("L68").Value = 7375
$ws.Range("M68").Value = -3750.5
$ws.Range("N68").Value = -8873
$ws.Range("H71").Value = 5937.25
$ws.Range("I71").Value = 4499.5
$ws.Range("J71").Value = 7375
$ws.Range("K71").Value = 22497.5
$ws.Range("L71").Value = 36875
$ws.Range("M71").Value = -18753.5
$ws.Range("N71").Value = -44363
$ws.Range("H122").Value = 835903.2
$ws.Range("I122").Value = 2604.125
$ws.Range("K122").Value = 7812.375
$ws.Range("M122").Value = -5362.375

# --- Sheet WVR: 25 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -858
$ws.Range("H126").Value = 4040
$ws.Range("I126").Value = 2983.3333
$ws.Range("J126").Value = 5625
$ws.Range("K126").Value = 8949.999899999999
$ws.Range("L126").Value = 16875
$ws.Range("M126").Value = -6479.999899999999
$ws.Range("N126").Value = -21815
$ws.Range("H132").Value = 28801.574
$ws.Range("I132").Value = 3209.3447
$ws.Range("J132").Value = 96272
$ws.Range("K132").Value = 9628.034100000001
$ws.Range("L132").Value = 288816
$ws.Range("M132").Value = -7098.034100000001
$ws.Range("N132").Value = -293876
$ws.Range("H136").Value = 9316623
$ws.Range("I136").Value = 11112389
$ws.Range("J136").Value = 337790.5
$ws.Range("K136").Value = 33337167
$ws.Range("L136").Value = 1013371.5
$ws.Range("M136").Value = -33334617
$ws.Range("N136").Value = -1018471.5
